$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the timestamp in A2
$ws.Range("A2").Value = "2025-05-27 09:03:41"

# Update the numeric metrics in row 2 (B2:AW2)
$ws.Range("B2").Value = 16232
$ws.Range("C2").Value = 11757
$ws.Range("D2").Value = 72.43100049285363
$ws.Range("E2").Value = 2243
$ws.Range("F2").Value = 13.81838344011828
$ws.Range("G2").Value = 3100
$ws.Range("H2").Value = 19.09807787087235
$ws.Range("I2").Value = 9560
$ws.Range("J2").Value = 58.89600788565797
$ws.Range("K2").Value = 3070327.87
$ws.Range("L2").Value = 3572
$ws.Range("M2").Value = 22.00591424346969
$ws.Range("N2").Value = 1171245.23
$ws.Range("O2").Value = 4698
$ws.Range("P2").Value = 28.942828979793
$ws.Range("Q2").Value = 499126.3
$ws.Range("R2").Value = 3572
$ws.Range("S2").Value = 22.00591424346969
$ws.Range("T2").Value = 3430
$ws.Range("U2").Value = 21.13109906357812
$ws.Range("V2").Value = 2416583.57
$ws.Range("W2").Value = 2300
$ws.Range("X2").Value = 14.1695416461311
$ws.Range("Y2").Value = 1432
$ws.Range("Z2").Value = 8.822079842286842
$ws.Range("AA2").Value = 154618
$ws.Range("AB2").Value = 800
$ws.Range("AC2").Value = 4.928536224741252
$ws.Range("AD2").Value = 16244
$ws.Range("AE2").Value = 10611
$ws.Range("AF2").Value = 65.32258064516128
$ws.Range("AG2").Value = 5633
$ws.Range("AH2").Value = 34.67741935483872
$ws.Range("AI2").Value = 484
$ws.Range("AJ2").Value = 863
$ws.Range("AK2").Value = 1348
$ws.Range("AL2").Value = 17.95918367346939
$ws.Range("AM2").Value = 32.02226345083488
$ws.Range("AN2").Value = 50.01855287569573
$ws.Range("AO2").Value = 1349074.39
$ws.Range("AP2").Value = 253575.28
$ws.Range("AQ2").Value = 84423.84
$ws.Range("AR2").Value = 79.96535906725249
$ws.Range("AS2").Value = 15.0304819853404
$ws.Range("AT2").Value = 5.004158947407098
$ws.Range("AU2").Value = 46.95067264573991
$ws.Range("AV2").Value = 230.1656394453005
$ws.Range("AW2").Value = 551.4683266202856
